$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3213.3333
$ws.Range("I40").Value = 2633.3333
$ws.Range("J40").Value = 3793.3333
$ws.Range("K40").Value = 2633.3333
$ws.Range("L40").Value = 3793.3333
$ws.Range("M40").Value = -2458.3333
$ws.Range("N40").Value = -4143.3333
$ws.Range("H98").Value = 727.1818
$ws.Range("I98").Value = 499.9
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 499.9
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 998.1
$ws.Range("N98").Value = -5996
$ws.Range("H111").Value = 6366.143
$ws.Range("I111").Value = 6658.25
$ws.Range("J111").Value = 5976.6665
$ws.Range("K111").Value = 19974.75
$ws.Range("L111").Value = 17929.9995
$ws.Range("M111").Value = -16907.75
$ws.Range("N111").Value = -24063.9995
$ws.Range("H122").Value = 727.1818
$ws.Range("I122").Value = 499.9
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 1499.7
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = 950.3000000000002
$ws.Range("N122").Value = -13900
$ws.Range("H125").Value = 1478.6666
$ws.Range("I125").Value = 783.3333
$ws.Range("J125").Value = 2174
$ws.Range("K125").Value = 7049.9997
$ws.Range("L125").Value = 19566
$ws.Range("M125").Value = -4589.9997
$ws.Range("N125").Value = -24486
$ws.Range("H129").Value = 962.90247
$ws.Range("I129").Value = 569.7273
$ws.Range("J129").Value = 1023.8169
$ws.Range("K129").Value = 1709.1819
$ws.Range("L129").Value = 3071.4507
$ws.Range("M129").Value = 3290.8181
$ws.Range("N129").Value = -13071.4507
$ws.Range("H135").Value = 715.9583
$ws.Range("I135").Value = 762.86365
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 6865.77285
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = -4330.77285
$ws.Range("N135").Value = -6870
$ws.Range("H137").Value = 3917.84
$ws.Range("I137").Value = 2652.65
$ws.Range("J137").Value = 8978.6
$ws.Range("K137").Value = 7957.950000000001
$ws.Range("L137").Value = 26935.8
$ws.Range("M137").Value = -5407.950000000001
$ws.Range("N137").Value = -32035.8
$ws.Range("H138").Value = 1861.4889
$ws.Range("J138").Value = 2130.647
$ws.Range("L138").Value = 6391.941
$ws.Range("N138").Value = -16671.941

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3088.6
$ws.Range("I2").Value = 3107.5
$ws.Range("K2").Value = 3107.5
$ws.Range("M2").Value = -2994.5
$ws.Range("H61").Value = 2451.913
$ws.Range("I61").Value = 1524.6875
$ws.Range("J61").Value = 4571.2856
$ws.Range("K61").Value = 1524.6875
$ws.Range("L61").Value = 4571.2856
$ws.Range("M61").Value = -1312.6875
$ws.Range("N61").Value = -4995.2856
$ws.Range("H102").Value = 1985.0714
$ws.Range("I102").Value = 1927.28
$ws.Range("J102").Value = 2466.6667
$ws.Range("K102").Value = 1927.28
$ws.Range("L102").Value = 2466.6667
$ws.Range("M102").Value = -305.28
$ws.Range("N102").Value = -5710.6667
$ws.Range("H116").Value = 3088.6
$ws.Range("I116").Value = 3107.5
$ws.Range("K116").Value = 3107.5
$ws.Range("M116").Value = -813.5
$ws.Range("H136").Value = 2451.913
$ws.Range("I136").Value = 1524.6875
$ws.Range("J136").Value = 4571.2856
$ws.Range("K136").Value = 4574.0625
$ws.Range("L136").Value = 13713.8568
$ws.Range("M136").Value = -2024.0625
$ws.Range("N136").Value = -18813.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3088.6
$ws.Range("I3").Value = 3107.5
$ws.Range("K3").Value = 3107.5
$ws.Range("M3").Value = -2993.5
$ws.Range("H57").Value = 11780
$ws.Range("J57").Value = 11780
$ws.Range("L57").Value = 11780
$ws.Range("N57").Value = -13220
$ws.Range("H107").Value = 1987.4117
$ws.Range("I107").Value = 1945.1538
$ws.Range("J107").Value = 2124.75
$ws.Range("K107").Value = 1945.1538
$ws.Range("L107").Value = 2124.75
$ws.Range("M107").Value = -25.15380000000005
$ws.Range("N107").Value = -5964.75
$ws.Range("H136").Value = 11780
$ws.Range("J136").Value = 11780
$ws.Range("L136").Value = 11780
$ws.Range("N136").Value = -21980

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6986.1704
$ws.Range("J31").Value = 10781.107
$ws.Range("L31").Value = 10781.107
$ws.Range("N31").Value = -11371.107
$ws.Range("H34").Value = 6986.1704
$ws.Range("J34").Value = 10781.107
$ws.Range("L34").Value = 10781.107
$ws.Range("N34").Value = -11185.107
$ws.Range("H58").Value = 1781.1666
$ws.Range("I58").Value = 1586
$ws.Range("J58").Value = 2757
$ws.Range("K58").Value = 1586
$ws.Range("L58").Value = 2757
$ws.Range("M58").Value = -1383
$ws.Range("N58").Value = -3163
$ws.Range("H94").Value = 1098.1875
$ws.Range("J94").Value = 1234.4546
$ws.Range("L94").Value = 1234.4546
$ws.Range("N94").Value = -2136.4546
$ws.Range("H134").Value = 2216.1155
$ws.Range("I134").Value = 1381.0952
$ws.Range("J134").Value = 5723.2
$ws.Range("K134").Value = 4143.2856
$ws.Range("L134").Value = 17169.6
$ws.Range("M134").Value = -1608.2856
$ws.Range("N134").Value = -22239.6
$ws.Range("H136").Value = 1781.1666
$ws.Range("I136").Value = 1586
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 4758
$ws.Range("L136").Value = 8271
$ws.Range("M136").Value = -2208
$ws.Range("N136").Value = -13371

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 355.875
$ws.Range("I5").Value = 355.875
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1067.625
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -955.625
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 6508.0586
$ws.Range("I122").Value = 545.6429000000001
$ws.Range("K122").Value = 4910.7861
$ws.Range("M122").Value = -2460.7861
$ws.Range("H132").Value = 1883.3158
$ws.Range("I132").Value = 1508.4286
$ws.Range("K132").Value = 13575.8574
$ws.Range("M132").Value = -11045.8574
$ws.Range("H135").Value = 355.875
$ws.Range("I135").Value = 355.875
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3202.875
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -667.875
$ws.Range("N135").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2607.2666
$ws.Range("I7").Value = 2238.3845
$ws.Range("K7").Value = 2238.3845
$ws.Range("M7").Value = -2126.3845
$ws.Range("H22").Value = 13231.9375
$ws.Range("I22").Value = 1540
$ws.Range("J22").Value = 18546.455
$ws.Range("K22").Value = 1540
$ws.Range("L22").Value = 18546.455
$ws.Range("M22").Value = -1245
$ws.Range("N22").Value = -19136.455
$ws.Range("H27").Value = 13231.9375
$ws.Range("I27").Value = 1540
$ws.Range("J27").Value = 18546.455
$ws.Range("K27").Value = 1540
$ws.Range("L27").Value = 18546.455
$ws.Range("M27").Value = -1433
$ws.Range("N27").Value = -18760.455
$ws.Range("H100").Value = 2668.5715
$ws.Range("I100").Value = 3653.3333
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 3653.3333
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -3112.3333
$ws.Range("N100").Value = -3482
$ws.Range("H122").Value = 2837.25
$ws.Range("I122").Value = 2426.8333
$ws.Range("J122").Value = 3576
$ws.Range("K122").Value = 7280.499899999999
$ws.Range("L122").Value = 10728
$ws.Range("M122").Value = -4830.499899999999
$ws.Range("N122").Value = -15628
$ws.Range("H126").Value = 2607.2666
$ws.Range("I126").Value = 2238.3845
$ws.Range("K126").Value = 6715.1535
$ws.Range("M126").Value = -4245.1535
$ws.Range("H137").Value = 38444.445
$ws.Range("I137").Value = 26000
$ws.Range("J137").Value = 40000
$ws.Range("K137").Value = 26000
$ws.Range("L137").Value = 40000
$ws.Range("M137").Value = -20900
$ws.Range("N137").Value = -50200

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 15142
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 15142
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15142
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -16124
$ws.Range("H113").Value = 915.9583
$ws.Range("I113").Value = 1136.7646
$ws.Range("J113").Value = 379.7143
$ws.Range("K113").Value = 3410.2938
$ws.Range("L113").Value = 1139.1429
$ws.Range("M113").Value = -1240.2938
$ws.Range("N113").Value = -5479.1429
$ws.Range("H122").Value = 2312.5356
$ws.Range("I122").Value = 1912.6666
$ws.Range("J122").Value = 2612.4375
$ws.Range("K122").Value = 5737.9998
$ws.Range("L122").Value = 7837.3125
$ws.Range("M122").Value = -3287.9998
$ws.Range("N122").Value = -12737.3125
$ws.Range("H132").Value = 4904306.5
$ws.Range("I132").Value = 2671.7778
$ws.Range("J132").Value = 10418646
$ws.Range("K132").Value = 8015.3334
$ws.Range("L132").Value = 31255938
$ws.Range("M132").Value = -5485.3334
$ws.Range("N132").Value = -31260998
$ws.Range("H136").Value = 2141.9019
$ws.Range("I136").Value = 1831.4054
$ws.Range("J136").Value = 2962.5
$ws.Range("K136").Value = 5494.216200000001
$ws.Range("L136").Value = 8887.5
$ws.Range("M136").Value = -2944.216200000001
$ws.Range("N136").Value = -13987.5
